# Update "想去人数" (want-to-go count) figures in column F across sheets.
# These values were refreshed by an automated scrape (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 1766
$ws.Range("F5").Value  = 445
$ws.Range("F8").Value  = 625
$ws.Range("F10").Value = 1705
$ws.Range("F11").Value = 345
$ws.Range("F13").Value = 798
$ws.Range("F16").Value = 12698
$ws.Range("F17").Value = 12718
$ws.Range("F18").Value = 945
$ws.Range("F19").Value = 736
$ws.Range("F20").Value = 9
$ws.Range("F21").Value = 503
$ws.Range("F23").Value = 525
$ws.Range("F27").Value = 238

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 52
$ws.Range("F10").Value = 71
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 6

# Sheet "全部类型" (All types - combined list)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 1766
$ws.Range("F7").Value  = 445
$ws.Range("F12").Value = 625
$ws.Range("F13").Value = 52
$ws.Range("F15").Value = 1705
$ws.Range("F16").Value = 345
$ws.Range("F18").Value = 798
$ws.Range("F22").Value = 12698
$ws.Range("F23").Value = 12718
$ws.Range("F24").Value = 945
$ws.Range("F25").Value = 736
$ws.Range("F26").Value = 9
$ws.Range("F27").Value = 503
$ws.Range("F29").Value = 525
$ws.Range("F37").Value = 238
$ws.Range("F39").Value = 71
$ws.Range("F40").Value = 5
$ws.Range("F41").Value = 6
